# "better answer form saving"
# Egor Barsukov's row (row 2) is re-scored (his answers to questions 1-7
# were actually blank, not ticked), and three more respondents' rows are
# appended (rows 4-6): Egor Barsukov again, Sergey Tsykura, and Fyodor
# Samokhin.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 2 (Егор Барсуков): columns C..I (questions 1-7) go from 1 to 0,
#     and the score totals in AD/AE drop accordingly.
for ($col = 3; $col -le 9; $col++) {
    $ws.Cells.Item(2, $col).Value = 0
}
$ws.Cells.Item(2, 30).Value = 5   # AD2 - primary score
$ws.Cells.Item(2, 31).Value = 34  # AE2 - secondary score

# --- Append new respondent rows 4, 5, 6 ---
function Set-ResponseRow {
    param($rowNum, $firstName, $lastName, $answers, $primary, $secondary)

    $ws.Cells.Item($rowNum, 1).Value = $firstName
    $ws.Cells.Item($rowNum, 2).Value = $lastName

    $col = 3
    foreach ($val in $answers) {
        $ws.Cells.Item($rowNum, $col).Value = $val
        $col = $col + 1
    }

    $ws.Cells.Item($rowNum, 30).Value = $primary
    $ws.Cells.Item($rowNum, 31).Value = $secondary
}

# Row 4: Егор Барсуков
Set-ResponseRow 4 "Егор" "Барсуков" `
    @(1,1,1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0) `
    4 27

# Row 5: Сергей Цыкура
Set-ResponseRow 5 "Сергей" "Цыкура" `
    @(0,1,1,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0) `
    4 27

# Row 6: Фёдор Самохин
Set-ResponseRow 6 "Фёдор" "Самохин" `
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0) `
    0 0
